$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) and volume-change (E) columns per the scraped data refresh.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.131.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.657.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5161"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("E7").Value = "  +0.25%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2628"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06263"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.650.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.423"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.00%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.883.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5409"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.34%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8102"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.71"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.152.69"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.004"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.613"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("E22").Value = "  -2.65%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.015"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.91"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1225"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.74%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.174"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.405"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.48%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05969"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.544"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.253"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.60%  "
$ws.Range("E34").Value = "  -5.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9641"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.84%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.425"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.769"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5675"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.14%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01592"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.964"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8544"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -7.99%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.798.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.64"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.980"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05172"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4183"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.19%  "
